$wb = $excel.ActiveWorkbook

# --- financial_statements sheet: update Example values (+.01) and insert new rows ---
$ws = $wb.Worksheets.Item("financial_statements")

$ws.Cells.Item(11, 1).Value = 'balance_sheet.current_assets.cash_and_equivalents'
$ws.Cells.Item(11, 2).Value = 'O montante total de caixa e equivalentes de caixa, incluindo moeda, contas bancárias e outros investimentos líquidos que podem ser rapidamente convertidos em dinheiro.'
$ws.Cells.Item(11, 3).Value = '''48572.01'
$ws.Cells.Item(11, 4).Value = 'number'
$ws.Cells.Item(11, 5).Value = 'float'
$ws.Cells.Item(11, 6).Value = 'Yes'
$ws.Cells.Item(11, 7).Value = 'Yes'

$ws.Cells.Item(12, 1).Value = 'balance_sheet.current_assets.short_term_investments'
$ws.Cells.Item(12, 2).Value = 'O valor dos investimentos que se espera serem liquidados em dinheiro dentro de um ano, como títulos negociáveis.'
$ws.Cells.Item(12, 3).Value = '''21345.01'
$ws.Cells.Item(12, 4).Value = 'number'
$ws.Cells.Item(12, 5).Value = 'float'
$ws.Cells.Item(12, 6).Value = 'Yes'
$ws.Cells.Item(12, 7).Value = 'Yes'

$ws.Cells.Item(13, 1).Value = 'balance_sheet.current_assets.accounts_receivable'
$ws.Cells.Item(13, 2).Value = 'O valor devido pelos clientes por vendas realizadas a crédito, que se espera receber em um curto período.'
$ws.Cells.Item(13, 3).Value = '''154321.01'
$ws.Cells.Item(13, 4).Value = 'number'
$ws.Cells.Item(13, 5).Value = 'float'
$ws.Cells.Item(13, 6).Value = 'Yes'
$ws.Cells.Item(13, 7).Value = 'Yes'

$ws.Cells.Item(14, 1).Value = 'balance_sheet.current_assets.notes_receivable'
$ws.Cells.Item(14, 2).Value = 'O valor de notas promissórias escritas recebidas de clientes ou outras partes, prometendo pagar um valor especificado até uma determinada data.'
$ws.Cells.Item(14, 3).Value = '''31789.01'
$ws.Cells.Item(14, 4).Value = 'number'
$ws.Cells.Item(14, 5).Value = 'float'
$ws.Cells.Item(14, 6).Value = 'Yes'
$ws.Cells.Item(14, 7).Value = 'Yes'

$ws.Cells.Item(15, 1).Value = 'balance_sheet.current_assets.other_debtors'
$ws.Cells.Item(15, 2).Value = 'Os valores totais devidos por vários outros devedores, excluindo contas a receber e notas promissórias.'
$ws.Cells.Item(15, 3).Value = '''12345.01'
$ws.Cells.Item(15, 4).Value = 'number'
$ws.Cells.Item(15, 5).Value = 'float'
$ws.Cells.Item(15, 6).Value = 'Yes'
$ws.Cells.Item(15, 7).Value = 'Yes'

$ws.Cells.Item(16, 1).Value = 'balance_sheet.current_assets.bad_debt_provision'
$ws.Cells.Item(16, 2).Value = 'O valor estimado de recebíveis que se espera serem incobráveis, frequentemente referido como provisão para devedores duvidosos.'
$ws.Cells.Item(16, 3).Value = '''0.01'
$ws.Cells.Item(16, 4).Value = 'number'
$ws.Cells.Item(16, 5).Value = 'float'
$ws.Cells.Item(16, 6).Value = 'Yes'
$ws.Cells.Item(16, 7).Value = 'Yes'

$ws.Cells.Item(17, 1).Value = 'balance_sheet.current_assets.tax_recoverable'
$ws.Cells.Item(17, 2).Value = 'O valor dos pagamentos de impostos que pode ser recuperado das autoridades fiscais.'
$ws.Cells.Item(17, 3).Value = '''8976.01'
$ws.Cells.Item(17, 4).Value = 'number'
$ws.Cells.Item(17, 5).Value = 'float'
$ws.Cells.Item(17, 6).Value = 'Yes'
$ws.Cells.Item(17, 7).Value = 'Yes'

$ws.Cells.Item(18, 1).Value = 'balance_sheet.current_assets.inventory'
$ws.Cells.Item(18, 2).Value = 'O valor total dos bens disponíveis para venda, matérias-primas, produtos em processo e produtos acabados.'
$ws.Cells.Item(18, 3).Value = '''65432.01'
$ws.Cells.Item(18, 4).Value = 'number'
$ws.Cells.Item(18, 5).Value = 'float'
$ws.Cells.Item(18, 6).Value = 'Yes'
$ws.Cells.Item(18, 7).Value = 'Yes'

$ws.Cells.Item(19, 1).Value = 'balance_sheet.current_assets.prepaid_expenses'
$ws.Cells.Item(19, 2).Value = 'O valor pago antecipadamente por bens ou serviços a serem recebidos no futuro, como prêmios de seguro ou aluguel.'
$ws.Cells.Item(19, 3).Value = '''14321.01'
$ws.Cells.Item(19, 4).Value = 'number'
$ws.Cells.Item(19, 5).Value = 'float'
$ws.Cells.Item(19, 6).Value = 'Yes'
$ws.Cells.Item(19, 7).Value = 'Yes'

$ws.Cells.Item(20, 1).Value = 'balance_sheet.current_assets.assets_available_for_sale'
$ws.Cells.Item(20, 2).Value = 'O valor dos ativos não circulantes que estão disponíveis para venda, mas ainda não foram vendidos, como equipamentos excedentes ou propriedades.'
$ws.Cells.Item(20, 3).Value = '''54321.01'
$ws.Cells.Item(20, 4).Value = 'number'
$ws.Cells.Item(20, 5).Value = 'float'
$ws.Cells.Item(20, 6).Value = 'Yes'
$ws.Cells.Item(20, 7).Value = 'Yes'

$ws.Cells.Item(21, 1).Value = 'balance_sheet.current_assets.total'
$ws.Cells.Item(21, 2).Value = 'A soma de todos os ativos circulantes, representando o valor total dos ativos que se espera serem convertidos em caixa ou utilizados dentro de um ano.'
$ws.Cells.Item(21, 3).Value = '''372480.01'
$ws.Cells.Item(21, 4).Value = 'number'
$ws.Cells.Item(21, 5).Value = 'float'
$ws.Cells.Item(21, 6).Value = 'Yes'
$ws.Cells.Item(21, 7).Value = 'Yes'

$ws.Cells.Item(22, 1).Value = 'balance_sheet.non_current_assets'
$ws.Cells.Item(22, 2).Value = 'Os ativos não circulantes da empresa, que são investimentos de longo prazo ou propriedades que não são facilmente convertidos em dinheiro, para o ano em questão.'
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = 'object'
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = 'No'
$ws.Cells.Item(22, 7).Value = 'No'

$ws.Cells.Item(23, 1).Value = 'balance_sheet.non_current_assets.property_plant_and_equipment'
$ws.Cells.Item(23, 2).Value = 'O valor total de propriedades, instalações e equipamentos de propriedade da empresa, incluindo terrenos, edifícios, maquinário e veículos, utilizados para operações de longo prazo.'
$ws.Cells.Item(23, 3).Value = '''1123456.01'
$ws.Cells.Item(23, 4).Value = 'number'
$ws.Cells.Item(23, 5).Value = 'float'
$ws.Cells.Item(23, 6).Value = 'Yes'
$ws.Cells.Item(23, 7).Value = 'Yes'

$ws.Cells.Item(24, 1).Value = 'balance_sheet.non_current_assets.accumulated_depreciation_and_amortization'
$ws.Cells.Item(24, 2).Value = 'Depreciação e amortização acumuladas totais, representando a alocação cumulativa do custo dos ativos não circulantes ao longo do período em que se espera que eles proporcionem benefícios econômicos.'
$ws.Cells.Item(24, 3).Value = '''123456.01'
$ws.Cells.Item(24, 4).Value = 'number'
$ws.Cells.Item(24, 5).Value = 'float'
$ws.Cells.Item(24, 6).Value = 'Yes'
$ws.Cells.Item(24, 7).Value = 'Yes'

$ws.Cells.Item(25, 1).Value = 'balance_sheet.non_current_assets.long_term_accounts_receivable'
$ws.Cells.Item(25, 2).Value = 'O valor devido pelos clientes por vendas feitas a crédito, esperado para ser recebido após um ano.'
$ws.Cells.Item(25, 3).Value = '''10987.01'
$ws.Cells.Item(25, 4).Value = 'number'
$ws.Cells.Item(25, 5).Value = 'float'
$ws.Cells.Item(25, 6).Value = 'Yes'
$ws.Cells.Item(25, 7).Value = 'Yes'

$ws.Cells.Item(26, 1).Value = 'balance_sheet.non_current_assets.prepayment_to_suppliers'
$ws.Cells.Item(26, 2).Value = 'O valor pago antecipadamente a fornecedores por bens ou serviços a serem recebidos no futuro, esperado para ser utilizado a longo prazo.'
$ws.Cells.Item(26, 3).Value = '''5432.01'
$ws.Cells.Item(26, 4).Value = 'number'
$ws.Cells.Item(26, 5).Value = 'float'
$ws.Cells.Item(26, 6).Value = 'Yes'
$ws.Cells.Item(26, 7).Value = 'Yes'

$ws.Cells.Item(27, 1).Value = 'balance_sheet.non_current_assets.goodwill'
$ws.Cells.Item(27, 2).Value = 'O valor dos ativos intangíveis que surgem da aquisição de outras empresas, representando o ágio pago sobre o valor justo dos ativos líquidos adquiridos.'
$ws.Cells.Item(27, 3).Value = '''47654.01'
$ws.Cells.Item(27, 4).Value = 'number'
$ws.Cells.Item(27, 5).Value = 'float'
$ws.Cells.Item(27, 6).Value = 'Yes'
$ws.Cells.Item(27, 7).Value = 'Yes'

$ws.Cells.Item(28, 1).Value = 'balance_sheet.non_current_assets.intangible_assets'
$ws.Cells.Item(28, 2).Value = 'O valor total dos ativos intangíveis de propriedade da empresa, como patentes, marcas registradas e direitos autorais, com vidas úteis que se estendem além de um ano.'
$ws.Cells.Item(28, 3).Value = '''43210.01'
$ws.Cells.Item(28, 4).Value = 'number'
$ws.Cells.Item(28, 5).Value = 'float'
$ws.Cells.Item(28, 6).Value = 'Yes'
$ws.Cells.Item(28, 7).Value = 'Yes'

$ws.Cells.Item(29, 1).Value = 'balance_sheet.non_current_assets.investments_in_associates'
$ws.Cells.Item(29, 2).Value = 'O valor dos investimentos em outras empresas nas quais a companhia tem influência significativa, mas não controle, tipicamente representado pela posse de 20-50% das ações com direito a voto do associado.'
$ws.Cells.Item(29, 3).Value = '''65432.01'
$ws.Cells.Item(29, 4).Value = 'number'
$ws.Cells.Item(29, 5).Value = 'float'
$ws.Cells.Item(29, 6).Value = 'Yes'
$ws.Cells.Item(29, 7).Value = 'Yes'

$ws.Cells.Item(30, 1).Value = 'balance_sheet.non_current_assets.long_term_financial_instruments'
$ws.Cells.Item(30, 2).Value = 'O valor dos instrumentos financeiros que se espera serem mantidos por mais de um ano, como títulos, debêntures e empréstimos de longo prazo.'
$ws.Cells.Item(30, 3).Value = '''32876.01'
$ws.Cells.Item(30, 4).Value = 'number'
$ws.Cells.Item(30, 5).Value = 'float'
$ws.Cells.Item(30, 6).Value = 'Yes'
$ws.Cells.Item(30, 7).Value = 'Yes'

$ws.Cells.Item(31, 1).Value = 'balance_sheet.non_current_assets.total'
$ws.Cells.Item(31, 2).Value = 'A soma de todos os ativos não circulantes, representando o valor total dos ativos que se espera serem utilizados ou mantidos por mais de um ano.'
$ws.Cells.Item(31, 3).Value = '''1346647.01'
$ws.Cells.Item(31, 4).Value = 'number'
$ws.Cells.Item(31, 5).Value = 'float'
$ws.Cells.Item(31, 6).Value = 'Yes'
$ws.Cells.Item(31, 7).Value = 'Yes'

$ws.Cells.Item(32, 1).Value = 'balance_sheet.current_liabilities'
$ws.Cells.Item(32, 2).Value = 'Os passivos circulantes da empresa, que se espera serem liquidados dentro do ano em questão.'
$ws.Cells.Item(32, 3).Value = ""
$ws.Cells.Item(32, 4).Value = 'object'
$ws.Cells.Item(32, 5).Value = ""
$ws.Cells.Item(32, 6).Value = 'No'
$ws.Cells.Item(32, 7).Value = 'No'

$ws.Cells.Item(33, 1).Value = 'balance_sheet.current_liabilities.bank_loans'
$ws.Cells.Item(33, 2).Value = 'O valor total dos empréstimos contraídos de bancos ou instituições financeiras, com expectativa de serem pagos dentro de um ano.'
$ws.Cells.Item(33, 3).Value = '''49876.01'
$ws.Cells.Item(33, 4).Value = 'number'
$ws.Cells.Item(33, 5).Value = 'float'
$ws.Cells.Item(33, 6).Value = 'Yes'
$ws.Cells.Item(33, 7).Value = 'Yes'

$ws.Cells.Item(34, 1).Value = 'balance_sheet.current_liabilities.accounts_payable'
$ws.Cells.Item(34, 2).Value = 'O valor devido aos fornecedores por bens ou serviços adquiridos a crédito, que se espera ser pago dentro de um curto período.'
$ws.Cells.Item(34, 3).Value = '''103298.01'
$ws.Cells.Item(34, 4).Value = 'number'
$ws.Cells.Item(34, 5).Value = 'float'
$ws.Cells.Item(34, 6).Value = 'Yes'
$ws.Cells.Item(34, 7).Value = 'Yes'

$ws.Cells.Item(35, 1).Value = 'balance_sheet.current_liabilities.notes_payable'
$ws.Cells.Item(35, 2).Value = 'O valor das notas promissórias escritas emitidas para fornecedores ou outros, prometendo pagar um valor especificado até uma determinada data.'
$ws.Cells.Item(35, 3).Value = '''25643.01'
$ws.Cells.Item(35, 4).Value = 'number'
$ws.Cells.Item(35, 5).Value = 'float'
$ws.Cells.Item(35, 6).Value = 'Yes'
$ws.Cells.Item(35, 7).Value = 'Yes'

$ws.Cells.Item(36, 1).Value = 'balance_sheet.current_liabilities.financial_instruments'
$ws.Cells.Item(36, 2).Value = 'O valor dos instrumentos financeiros que se espera serem liquidados em dinheiro dentro de um ano, como títulos, debêntures e empréstimos de curto prazo.'
$ws.Cells.Item(36, 3).Value = '''14321.01'
$ws.Cells.Item(36, 4).Value = 'number'
$ws.Cells.Item(36, 5).Value = 'float'
$ws.Cells.Item(36, 6).Value = 'Yes'
$ws.Cells.Item(36, 7).Value = 'Yes'

$ws.Cells.Item(37, 1).Value = 'balance_sheet.current_liabilities.other_creditors'
$ws.Cells.Item(37, 2).Value = 'Os valores totais devidos a vários outros credores, excluindo contas a pagar e notas promissórias.'
$ws.Cells.Item(37, 3).Value = '''21987.01'
$ws.Cells.Item(37, 4).Value = 'number'
$ws.Cells.Item(37, 5).Value = 'float'
$ws.Cells.Item(37, 6).Value = 'Yes'
$ws.Cells.Item(37, 7).Value = 'Yes'

$ws.Cells.Item(38, 1).Value = 'balance_sheet.current_liabilities.income_tax_payable'
$ws.Cells.Item(38, 2).Value = 'O valor do imposto de renda devido às autoridades fiscais, que se espera ser pago em um curto período.'
$ws.Cells.Item(38, 3).Value = '''12765.01'
$ws.Cells.Item(38, 4).Value = 'number'
$ws.Cells.Item(38, 5).Value = 'float'
$ws.Cells.Item(38, 6).Value = 'Yes'
$ws.Cells.Item(38, 7).Value = 'Yes'

$ws.Cells.Item(39, 1).Value = 'balance_sheet.current_liabilities.customer_advances'
$ws.Cells.Item(39, 2).Value = 'O valor total recebido antecipadamente de clientes por bens ou serviços a serem entregues no futuro, esperado para ser utilizado dentro de um ano.'
$ws.Cells.Item(39, 3).Value = '''18765.01'
$ws.Cells.Item(39, 4).Value = 'number'
$ws.Cells.Item(39, 5).Value = 'float'
$ws.Cells.Item(39, 6).Value = 'Yes'
$ws.Cells.Item(39, 7).Value = 'Yes'

$ws.Cells.Item(40, 1).Value = 'balance_sheet.current_liabilities.provisions'
$ws.Cells.Item(40, 2).Value = 'O valor estimado reservado para passivos ou perdas futuras, como garantias, reivindicações legais ou custos de reestruturação.'
$ws.Cells.Item(40, 3).Value = '''10987.01'
$ws.Cells.Item(40, 4).Value = 'number'
$ws.Cells.Item(40, 5).Value = 'float'
$ws.Cells.Item(40, 6).Value = 'Yes'
$ws.Cells.Item(40, 7).Value = 'Yes'

$ws.Cells.Item(41, 1).Value = 'balance_sheet.current_liabilities.taxes_payable'
$ws.Cells.Item(41, 2).Value = 'O valor total de impostos devidos às autoridades fiscais, que se espera ser pago em um curto período.'
$ws.Cells.Item(41, 3).Value = '''5321.01'
$ws.Cells.Item(41, 4).Value = 'number'
$ws.Cells.Item(41, 5).Value = 'float'
$ws.Cells.Item(41, 6).Value = 'Yes'
$ws.Cells.Item(41, 7).Value = 'Yes'

$ws.Cells.Item(42, 1).Value = 'balance_sheet.current_liabilities.total'
$ws.Cells.Item(42, 2).Value = 'A soma de todos os passivos circulantes, representando o valor total das obrigações que se espera serem liquidadas dentro de um ano.'
$ws.Cells.Item(42, 3).Value = '''260963.01'
$ws.Cells.Item(42, 4).Value = 'number'
$ws.Cells.Item(42, 5).Value = 'float'
$ws.Cells.Item(42, 6).Value = 'Yes'
$ws.Cells.Item(42, 7).Value = 'Yes'

$ws.Cells.Item(43, 1).Value = 'balance_sheet.non_current_liabilities'
$ws.Cells.Item(43, 2).Value = 'As obrigações não circulantes da empresa, que são obrigações de longo prazo não exigíveis dentro do ano em questão.'
$ws.Cells.Item(43, 3).Value = ""
$ws.Cells.Item(43, 4).Value = 'object'
$ws.Cells.Item(43, 5).Value = ""
$ws.Cells.Item(43, 6).Value = 'No'
$ws.Cells.Item(43, 7).Value = 'No'

$ws.Cells.Item(44, 1).Value = 'balance_sheet.non_current_liabilities.long_term_accounts_payable'
$ws.Cells.Item(44, 2).Value = 'O valor devido aos fornecedores por bens ou serviços adquiridos a crédito, com expectativa de pagamento após um ano.'
$ws.Cells.Item(44, 3).Value = '''30876.01'
$ws.Cells.Item(44, 4).Value = 'number'
$ws.Cells.Item(44, 5).Value = 'float'
$ws.Cells.Item(44, 6).Value = 'Yes'
$ws.Cells.Item(44, 7).Value = 'Yes'

$ws.Cells.Item(45, 1).Value = 'balance_sheet.non_current_liabilities.long_term_financial_instruments'
$ws.Cells.Item(45, 2).Value = 'O valor dos instrumentos financeiros que se espera serem mantidos por mais de um ano, como títulos, debêntures e empréstimos de longo prazo.'
$ws.Cells.Item(45, 3).Value = '''42310.01'
$ws.Cells.Item(45, 4).Value = 'number'
$ws.Cells.Item(45, 5).Value = 'float'
$ws.Cells.Item(45, 6).Value = 'Yes'
$ws.Cells.Item(45, 7).Value = 'Yes'

$ws.Cells.Item(46, 1).Value = 'balance_sheet.non_current_liabilities.deferred_revenue'
$ws.Cells.Item(46, 2).Value = 'O montante recebido antecipadamente de clientes por bens ou serviços a serem entregues no futuro, esperado para ser reconhecido como receita a longo prazo (como aluguel).'
$ws.Cells.Item(46, 3).Value = '''21987.01'
$ws.Cells.Item(46, 4).Value = 'number'
$ws.Cells.Item(46, 5).Value = 'float'
$ws.Cells.Item(46, 6).Value = 'Yes'
$ws.Cells.Item(46, 7).Value = 'Yes'

$ws.Cells.Item(47, 1).Value = 'balance_sheet.non_current_liabilities.contributions_for_future_capital_increases'
$ws.Cells.Item(47, 2).Value = 'As contribuições totais recebidas de acionistas ou outros investidores para futuros aumentos de capital, que se espera serem utilizadas a longo prazo.'
$ws.Cells.Item(47, 3).Value = '''10987.01'
$ws.Cells.Item(47, 4).Value = 'number'
$ws.Cells.Item(47, 5).Value = 'float'
$ws.Cells.Item(47, 6).Value = 'Yes'
$ws.Cells.Item(47, 7).Value = 'Yes'

$ws.Cells.Item(48, 1).Value = 'balance_sheet.non_current_liabilities.deferred_income_tax'
$ws.Cells.Item(48, 2).Value = 'O valor do imposto de renda que é diferido para períodos futuros, esperado para ser pago após um ano.'
$ws.Cells.Item(48, 3).Value = '''26543.01'
$ws.Cells.Item(48, 4).Value = 'number'
$ws.Cells.Item(48, 5).Value = 'float'
$ws.Cells.Item(48, 6).Value = 'Yes'
$ws.Cells.Item(48, 7).Value = 'Yes'

$ws.Cells.Item(49, 1).Value = 'balance_sheet.non_current_liabilities.employee_benefits'
$ws.Cells.Item(49, 2).Value = 'O valor total dos benefícios devidos aos empregados, como pensões, gratificações e outros benefícios pós-emprego, que se espera serem liquidados a longo prazo.'
$ws.Cells.Item(49, 3).Value = '''30218.01'
$ws.Cells.Item(49, 4).Value = 'number'
$ws.Cells.Item(49, 5).Value = 'float'
$ws.Cells.Item(49, 6).Value = 'Yes'
$ws.Cells.Item(49, 7).Value = 'Yes'

$ws.Cells.Item(50, 1).Value = 'balance_sheet.non_current_liabilities.long_term_provisions'
$ws.Cells.Item(50, 2).Value = 'O valor estimado reservado para passivos ou perdas futuras, como garantias, reivindicações legais ou custos de reestruturação, que se espera serem liquidados após um ano.'
$ws.Cells.Item(50, 3).Value = '''15432.01'
$ws.Cells.Item(50, 4).Value = 'number'
$ws.Cells.Item(50, 5).Value = 'float'
$ws.Cells.Item(50, 6).Value = 'Yes'
$ws.Cells.Item(50, 7).Value = 'Yes'

$ws.Cells.Item(51, 1).Value = 'balance_sheet.non_current_liabilities.total'
$ws.Cells.Item(51, 2).Value = 'A soma de todos os passivos não circulantes, representando o valor total das obrigações que se espera serem liquidadas após um ano.'
$ws.Cells.Item(51, 3).Value = '''178353.01'
$ws.Cells.Item(51, 4).Value = 'number'
$ws.Cells.Item(51, 5).Value = 'float'
$ws.Cells.Item(51, 6).Value = 'Yes'
$ws.Cells.Item(51, 7).Value = 'Yes'

$ws.Cells.Item(52, 1).Value = 'balance_sheet.equity'
$ws.Cells.Item(52, 2).Value = 'O patrimônio líquido da empresa, representando o interesse residual nos ativos após a dedução dos passivos.'
$ws.Cells.Item(52, 3).Value = ""
$ws.Cells.Item(52, 4).Value = 'object'
$ws.Cells.Item(52, 5).Value = ""
$ws.Cells.Item(52, 6).Value = 'No'
$ws.Cells.Item(52, 7).Value = 'No'

$ws.Cells.Item(53, 1).Value = 'balance_sheet.equity.stockholders_equity'
$ws.Cells.Item(53, 2).Value = 'O valor total das ações emitidas pela empresa, representando o interesse de propriedade dos acionistas no negócio.'
$ws.Cells.Item(53, 3).Value = '''501234.01'
$ws.Cells.Item(53, 4).Value = 'number'
$ws.Cells.Item(53, 5).Value = 'float'
$ws.Cells.Item(53, 6).Value = 'Yes'
$ws.Cells.Item(53, 7).Value = 'Yes'

$ws.Cells.Item(54, 1).Value = 'balance_sheet.equity.future_capital_contributions'
$ws.Cells.Item(54, 2).Value = 'Os fundos recebidos dos acionistas que são especificamente designados para futuros aumentos de capital ou investimentos.'
$ws.Cells.Item(54, 3).Value = '''75000.01'
$ws.Cells.Item(54, 4).Value = 'number'
$ws.Cells.Item(54, 5).Value = 'float'
$ws.Cells.Item(54, 6).Value = 'Yes'
$ws.Cells.Item(54, 7).Value = 'Yes'

$ws.Cells.Item(55, 1).Value = 'balance_sheet.equity.legal_reserve'
$ws.Cells.Item(55, 2).Value = 'A reserva legal exigida por lei, geralmente separada dos lucros, para fornecer proteção financeira contra perdas ou obrigações futuras.'
$ws.Cells.Item(55, 3).Value = '''25000.01'
$ws.Cells.Item(55, 4).Value = 'number'
$ws.Cells.Item(55, 5).Value = 'float'
$ws.Cells.Item(55, 6).Value = 'Yes'
$ws.Cells.Item(55, 7).Value = 'Yes'

$ws.Cells.Item(56, 1).Value = 'balance_sheet.equity.capital_update_excess'
$ws.Cells.Item(56, 2).Value = 'O excedente resultante de ajustes feitos no capital próprio, frequentemente devido à inflação ou à reavaliação de ativos.'
$ws.Cells.Item(56, 3).Value = '''15000.01'
$ws.Cells.Item(56, 4).Value = 'number'
$ws.Cells.Item(56, 5).Value = 'float'
$ws.Cells.Item(56, 6).Value = 'Yes'
$ws.Cells.Item(56, 7).Value = 'Yes'

$ws.Cells.Item(57, 1).Value = 'balance_sheet.equity.capital_update_insufficiency'
$ws.Cells.Item(57, 2).Value = 'O déficit resultante de ajustes feitos no capital próprio, muitas vezes devido à inflação ou à reavaliação de ativos.'
$ws.Cells.Item(57, 3).Value = '''-5000.01'
$ws.Cells.Item(57, 4).Value = 'number'
$ws.Cells.Item(57, 5).Value = 'float'
$ws.Cells.Item(57, 6).Value = 'Yes'
$ws.Cells.Item(57, 7).Value = 'Yes'

$ws.Cells.Item(58, 1).Value = 'balance_sheet.equity.capital_reserve'
$ws.Cells.Item(58, 2).Value = 'A reserva de capital derivada de atividades não operacionais, como ganhos de reavaliações de ativos ou certas transações de capital.'
$ws.Cells.Item(58, 3).Value = '''10000.01'
$ws.Cells.Item(58, 4).Value = 'number'
$ws.Cells.Item(58, 5).Value = 'float'
$ws.Cells.Item(58, 6).Value = 'Yes'
$ws.Cells.Item(58, 7).Value = 'Yes'

$ws.Cells.Item(59, 1).Value = 'balance_sheet.equity.share_premium_on_stock_sales'
$ws.Cells.Item(59, 2).Value = 'O valor excedente recebido por uma empresa quando as ações são emitidas a um preço acima do seu valor nominal (par).'
$ws.Cells.Item(59, 3).Value = '''50000.01'
$ws.Cells.Item(59, 4).Value = 'number'
$ws.Cells.Item(59, 5).Value = 'float'
$ws.Cells.Item(59, 6).Value = 'Yes'
$ws.Cells.Item(59, 7).Value = 'Yes'

$ws.Cells.Item(60, 1).Value = 'balance_sheet.equity.retained_earnings'
$ws.Cells.Item(60, 2).Value = 'Os lucros ou prejuízos acumulados da empresa que não foram distribuídos aos acionistas como dividendos.'
$ws.Cells.Item(60, 3).Value = '''202345.01'
$ws.Cells.Item(60, 4).Value = 'number'
$ws.Cells.Item(60, 5).Value = 'float'
$ws.Cells.Item(60, 6).Value = 'Yes'
$ws.Cells.Item(60, 7).Value = 'Yes'

$ws.Cells.Item(61, 1).Value = 'balance_sheet.equity.other_comprehensive_income'
$ws.Cells.Item(61, 2).Value = 'Os ganhos ou perdas que não estão incluídos no lucro líquido, mas são reportados diretamente no patrimônio, como ganhos não realizados em investimentos ou ajustes de tradução de moeda estrangeira.'
$ws.Cells.Item(61, 3).Value = '''10987.01'
$ws.Cells.Item(61, 4).Value = 'number'
$ws.Cells.Item(61, 5).Value = 'float'
$ws.Cells.Item(61, 6).Value = 'Yes'
$ws.Cells.Item(61, 7).Value = 'Yes'

$ws.Cells.Item(62, 1).Value = 'balance_sheet.equity.controlling_interest'
$ws.Cells.Item(62, 2).Value = 'O interesse de propriedade na empresa detido pela entidade controladora ou pelos acionistas majoritários, representando a participação controladora no negócio.'
$ws.Cells.Item(62, 3).Value = '''70876.01'
$ws.Cells.Item(62, 4).Value = 'number'
$ws.Cells.Item(62, 5).Value = 'float'
$ws.Cells.Item(62, 6).Value = 'Yes'
$ws.Cells.Item(62, 7).Value = 'Yes'

$ws.Cells.Item(63, 1).Value = 'balance_sheet.equity.non_controlling_interest'
$ws.Cells.Item(63, 2).Value = 'O interesse de propriedade na empresa detido por acionistas minoritários, representando a participação não controladora no negócio.'
$ws.Cells.Item(63, 3).Value = '''50321.01'
$ws.Cells.Item(63, 4).Value = 'number'
$ws.Cells.Item(63, 5).Value = 'float'
$ws.Cells.Item(63, 6).Value = 'Yes'
$ws.Cells.Item(63, 7).Value = 'Yes'

$ws.Cells.Item(64, 1).Value = 'balance_sheet.equity.total'
$ws.Cells.Item(64, 2).Value = 'A soma do capital social, lucros retidos, outros resultados abrangentes, participação controladora e participação não controladora, representando o patrimônio total da empresa.'
$ws.Cells.Item(64, 3).Value = '''836763.01'
$ws.Cells.Item(64, 4).Value = 'number'
$ws.Cells.Item(64, 5).Value = 'float'
$ws.Cells.Item(64, 6).Value = 'Yes'
$ws.Cells.Item(64, 7).Value = 'Yes'

$ws.Cells.Item(65, 1).Value = 'income_statement'
$ws.Cells.Item(65, 2).Value = 'A demonstração de resultados detalhando as receitas, despesas e lucros da empresa para o ano em questão.'
$ws.Cells.Item(65, 3).Value = ""
$ws.Cells.Item(65, 4).Value = 'object'
$ws.Cells.Item(65, 5).Value = ""
$ws.Cells.Item(65, 6).Value = 'Yes'
$ws.Cells.Item(65, 7).Value = 'No'

$ws.Cells.Item(66, 1).Value = 'income_statement.net_revenue'
$ws.Cells.Item(66, 2).Value = 'A receita total gerada pela empresa a partir de suas operações principais, excluindo quaisquer deduções para descontos, devoluções ou abatimentos.

> **Nota**: `domestic_sales` + `foreign_sales` não somarão o `net_revenue` devido à exclusão de descontos, devoluções e abatimentos.'
$ws.Cells.Item(66, 3).Value = '''1212345.01'
$ws.Cells.Item(66, 4).Value = 'number'
$ws.Cells.Item(66, 5).Value = 'float'
$ws.Cells.Item(66, 6).Value = 'Yes'
$ws.Cells.Item(66, 7).Value = 'Yes'

$ws.Cells.Item(67, 1).Value = 'income_statement.domestic_sales'
$ws.Cells.Item(67, 2).Value = 'A receita gerada pela empresa a partir da venda de bens ou serviços dentro de seu país de origem.'
$ws.Cells.Item(67, 3).Value = '''1123456.01'
$ws.Cells.Item(67, 4).Value = 'number'
$ws.Cells.Item(67, 5).Value = 'float'
$ws.Cells.Item(67, 6).Value = 'Yes'
$ws.Cells.Item(67, 7).Value = 'Yes'

$ws.Cells.Item(68, 1).Value = 'income_statement.foreign_sales'
$ws.Cells.Item(68, 2).Value = 'A receita gerada pela empresa com a venda de bens ou serviços em países estrangeiros.'
$ws.Cells.Item(68, 3).Value = '''88987.01'
$ws.Cells.Item(68, 4).Value = 'number'
$ws.Cells.Item(68, 5).Value = 'float'
$ws.Cells.Item(68, 6).Value = 'Yes'
$ws.Cells.Item(68, 7).Value = 'Yes'

$ws.Cells.Item(69, 1).Value = 'income_statement.materials_used'
$ws.Cells.Item(69, 2).Value = 'O custo total dos materiais utilizados ou comercializados pela empresa durante o período de relatório.'
$ws.Cells.Item(69, 3).Value = '''609876.01'
$ws.Cells.Item(69, 4).Value = 'number'
$ws.Cells.Item(69, 5).Value = 'float'
$ws.Cells.Item(69, 6).Value = 'No'
$ws.Cells.Item(69, 7).Value = 'Yes'

$ws.Cells.Item(70, 1).Value = 'income_statement.cost_of_goods_sold'
$ws.Cells.Item(70, 2).Value = 'O custo total incorrido pela empresa para produzir ou adquirir os bens vendidos durante o período de relatório.'
$ws.Cells.Item(70, 3).Value = '''412345.01'
$ws.Cells.Item(70, 4).Value = 'number'
$ws.Cells.Item(70, 5).Value = 'float'
$ws.Cells.Item(70, 6).Value = 'Yes'
$ws.Cells.Item(70, 7).Value = 'Yes'

$ws.Cells.Item(71, 1).Value = 'income_statement.cost_of_services_sold'
$ws.Cells.Item(71, 2).Value = 'O custo total incorrido pela empresa para fornecer os serviços vendidos durante o período de relatório.'
$ws.Cells.Item(71, 3).Value = '''101234.01'
$ws.Cells.Item(71, 4).Value = 'number'
$ws.Cells.Item(71, 5).Value = 'float'
$ws.Cells.Item(71, 6).Value = 'Yes'
$ws.Cells.Item(71, 7).Value = 'Yes'

$ws.Cells.Item(72, 1).Value = 'income_statement.gross_profit'
$ws.Cells.Item(72, 2).Value = 'A diferença entre a receita líquida e o custo total dos bens e serviços vendidos, representando o lucro obtido com as operações principais do negócio antes de deduzir as despesas operacionais.'
$ws.Cells.Item(72, 3).Value = '''190890.01'
$ws.Cells.Item(72, 4).Value = 'number'
$ws.Cells.Item(72, 5).Value = 'float'
$ws.Cells.Item(72, 6).Value = 'Yes'
$ws.Cells.Item(72, 7).Value = 'Yes'

$ws.Cells.Item(73, 1).Value = 'income_statement.gross_loss'
$ws.Cells.Item(73, 2).Value = 'A diferença negativa entre a receita líquida e o custo total dos bens e serviços vendidos, representando a perda incorrida nas operações principais do negócio antes de deduzir as despesas operacionais.'
$ws.Cells.Item(73, 3).Value = ""
$ws.Cells.Item(73, 4).Value = 'number'
$ws.Cells.Item(73, 5).Value = 'float'
$ws.Cells.Item(73, 6).Value = 'Yes'
$ws.Cells.Item(73, 7).Value = 'Yes'

$ws.Cells.Item(74, 1).Value = 'income_statement.operating_expenses'
$ws.Cells.Item(74, 2).Value = 'As despesas totais incorridas pela empresa em suas atividades operacionais normais, incluindo despesas de vendas, gerais e administrativas.'
$ws.Cells.Item(74, 3).Value = '''122345.01'
$ws.Cells.Item(74, 4).Value = 'number'
$ws.Cells.Item(74, 5).Value = 'float'
$ws.Cells.Item(74, 6).Value = 'Yes'
$ws.Cells.Item(74, 7).Value = 'Yes'

$ws.Cells.Item(75, 1).Value = 'income_statement.operating_income'
$ws.Cells.Item(75, 2).Value = 'O lucro obtido das operações principais do negócio após a dedução das despesas operacionais, mas antes de considerar juros, impostos e outros itens não operacionais.'
$ws.Cells.Item(75, 3).Value = '''68545.01'
$ws.Cells.Item(75, 4).Value = 'number'
$ws.Cells.Item(75, 5).Value = 'float'
$ws.Cells.Item(75, 6).Value = 'Yes'
$ws.Cells.Item(75, 7).Value = 'Yes'

$ws.Cells.Item(76, 1).Value = 'income_statement.operating_loss'
$ws.Cells.Item(76, 2).Value = 'O prejuízo incorrido nas operações principais do negócio após deduzir as despesas operacionais, mas antes de considerar juros, impostos e outros itens não operacionais.'
$ws.Cells.Item(76, 3).Value = ""
$ws.Cells.Item(76, 4).Value = 'number'
$ws.Cells.Item(76, 5).Value = 'float'
$ws.Cells.Item(76, 6).Value = 'Yes'
$ws.Cells.Item(76, 7).Value = 'Yes'

$ws.Cells.Item(77, 1).Value = 'income_statement.financial_result'
$ws.Cells.Item(77, 2).Value = 'O resultado líquido das atividades financeiras, incluindo receita de juros, despesa de juros e outros ganhos ou perdas financeiras.'
$ws.Cells.Item(77, 3).Value = '''15098.01'
$ws.Cells.Item(77, 4).Value = 'number'
$ws.Cells.Item(77, 5).Value = 'float'
$ws.Cells.Item(77, 6).Value = 'Yes'
$ws.Cells.Item(77, 7).Value = 'Yes'

$ws.Cells.Item(78, 1).Value = 'income_statement.income_statement_financial_gains'
$ws.Cells.Item(78, 2).Value = 'A receita financeira total positiva, incluindo receita de juros, ganhos cambiais e outros ganhos de atividades de financiamento. Este valor deve ser sempre positivo.'
$ws.Cells.Item(78, 3).Value = '''85000.01'
$ws.Cells.Item(78, 4).Value = 'number'
$ws.Cells.Item(78, 5).Value = 'float'
$ws.Cells.Item(78, 6).Value = 'Yes'
$ws.Cells.Item(78, 7).Value = 'Yes'

$ws.Cells.Item(79, 1).Value = 'income_statement.income_statement_financial_costs'
$ws.Cells.Item(79, 2).Value = 'As despesas financeiras totais, incluindo despesas com juros, perdas cambiais e outros custos incorridos em atividades de financiamento. Este valor deve ser sempre negativo.'
$ws.Cells.Item(79, 3).Value = '''-32000.01'
$ws.Cells.Item(79, 4).Value = 'number'
$ws.Cells.Item(79, 5).Value = 'float'
$ws.Cells.Item(79, 6).Value = 'Yes'
$ws.Cells.Item(79, 7).Value = 'Yes'

$ws.Cells.Item(80, 1).Value = 'income_statement.equity_in_earnings_of_affiliates'
$ws.Cells.Item(80, 2).Value = 'A participação da empresa no lucro ou prejuízo de suas associadas, entidades sobre as quais ela tem influência significativa, mas não controle.'
$ws.Cells.Item(80, 3).Value = '''5678.01'
$ws.Cells.Item(80, 4).Value = 'number'
$ws.Cells.Item(80, 5).Value = 'float'
$ws.Cells.Item(80, 6).Value = 'Yes'
$ws.Cells.Item(80, 7).Value = 'Yes'

$ws.Cells.Item(81, 1).Value = 'income_statement.income_before_taxes'
$ws.Cells.Item(81, 2).Value = 'O lucro obtido antes de contabilizar as despesas com imposto de renda.'
$ws.Cells.Item(81, 3).Value = '''89321.01'
$ws.Cells.Item(81, 4).Value = 'number'
$ws.Cells.Item(81, 5).Value = 'float'
$ws.Cells.Item(81, 6).Value = 'Yes'
$ws.Cells.Item(81, 7).Value = 'Yes'

$ws.Cells.Item(82, 1).Value = 'income_statement.loss_before_taxes'
$ws.Cells.Item(82, 2).Value = 'O prejuízo incorrido antes de contabilizar as despesas com imposto de renda.'
$ws.Cells.Item(82, 3).Value = ""
$ws.Cells.Item(82, 4).Value = 'number'
$ws.Cells.Item(82, 5).Value = 'float'
$ws.Cells.Item(82, 6).Value = 'Yes'
$ws.Cells.Item(82, 7).Value = 'Yes'

$ws.Cells.Item(83, 1).Value = 'income_statement.income_taxes'
$ws.Cells.Item(83, 2).Value = 'O valor total das despesas com imposto de renda incorridas durante o período de relatório.'
$ws.Cells.Item(83, 3).Value = '''20123.01'
$ws.Cells.Item(83, 4).Value = 'number'
$ws.Cells.Item(83, 5).Value = 'float'
$ws.Cells.Item(83, 6).Value = 'Yes'
$ws.Cells.Item(83, 7).Value = 'Yes'

$ws.Cells.Item(84, 1).Value = 'income_statement.income_from_continuing_operations'
$ws.Cells.Item(84, 2).Value = 'O lucro obtido das operações comerciais contínuas da empresa após a dedução das despesas operacionais e impostos.'
$ws.Cells.Item(84, 3).Value = '''69198.01'
$ws.Cells.Item(84, 4).Value = 'number'
$ws.Cells.Item(84, 5).Value = 'float'
$ws.Cells.Item(84, 6).Value = 'Yes'
$ws.Cells.Item(84, 7).Value = 'Yes'

$ws.Cells.Item(85, 1).Value = 'income_statement.loss_from_continuing_operations'
$ws.Cells.Item(85, 2).Value = 'O prejuízo incorrido pelas operações comerciais contínuas da empresa após a dedução das despesas operacionais e impostos.'
$ws.Cells.Item(85, 3).Value = ""
$ws.Cells.Item(85, 4).Value = 'number'
$ws.Cells.Item(85, 5).Value = 'float'
$ws.Cells.Item(85, 6).Value = 'Yes'
$ws.Cells.Item(85, 7).Value = 'Yes'

$ws.Cells.Item(86, 1).Value = 'income_statement.discontinued_operations'
$ws.Cells.Item(86, 2).Value = 'O resultado líquido das operações que foram descontinuadas ou vendidas durante o período de relatório.'
$ws.Cells.Item(86, 3).Value = '''0.01'
$ws.Cells.Item(86, 4).Value = 'number'
$ws.Cells.Item(86, 5).Value = 'float'
$ws.Cells.Item(86, 6).Value = 'Yes'
$ws.Cells.Item(86, 7).Value = 'Yes'

$ws.Cells.Item(87, 1).Value = 'income_statement.net_income'
$ws.Cells.Item(87, 2).Value = 'O lucro total obtido pela empresa após deduzir todas as despesas, incluindo operacionais, não operacionais, juros e impostos.'
$ws.Cells.Item(87, 3).Value = '''69198.01'
$ws.Cells.Item(87, 4).Value = 'number'
$ws.Cells.Item(87, 5).Value = 'float'
$ws.Cells.Item(87, 6).Value = 'Yes'
$ws.Cells.Item(87, 7).Value = 'Yes'

$ws.Cells.Item(88, 1).Value = 'income_statement.net_loss'
$ws.Cells.Item(88, 2).Value = 'O prejuízo total incorrido pela empresa após deduzir todas as despesas, incluindo operacionais, não operacionais, juros e impostos.'
$ws.Cells.Item(88, 3).Value = ""
$ws.Cells.Item(88, 4).Value = 'number'
$ws.Cells.Item(88, 5).Value = 'float'
$ws.Cells.Item(88, 6).Value = 'Yes'
$ws.Cells.Item(88, 7).Value = 'Yes'

# --- links sheet: tweak credentials_storage description wording ---
$wsLinks = $wb.Worksheets.Item("links")
$wsLinks.Range("B12").Value = 'Indica se as credenciais devem ou não ser armazenadas (e a duração para a qual as credenciais serão armazenadas).

- Para links recorrentes, isso é definido como `store` por padrão (e não pode ser alterado).
- Para links únicos, isso é definido como `365d` por padrão.

Pode ser:
  - `store` para armazenar credenciais (até que o link seja excluído)
  - `nostore` para não armazenar credenciais
  - Qualquer valor entre `1d` e `365d` para indicar o número de dias que você deseja que as credenciais sejam armazenadas.

Para mais informações, confira a seção <a href="https://developers.belvo.com/docs/data-retention-controls#credentials_storage" target="_blank">credentials_storage</a> do nosso artigo sobre controles de retenção de dados.'
